# Rows 34-50 on the "Artfynd" sheet had their per-record data shuffled among
# themselves (row 44 keeps its own data; every other row now shows the data
# that used to belong to a different row in that same block). Columns that
# are identical for every record in this block (e.g. C, P, S, T, U, V, W, Y,
# AA, AD, AE, AG) or that are always blank (J, K, L, N, O, X, AC, AF, AH-AS,
# AU, AV) or always an explicit empty string (I, AT, AY) are left completely
# untouched so their exact original representation is preserved. Only the
# columns whose values actually differ from record to record are rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (whose record now lives in destRow).
$mapping = @{
    34 = 45
    35 = 46
    36 = 47
    37 = 35
    38 = 40
    39 = 50
    40 = 36
    41 = 37
    42 = 34
    43 = 39
    44 = 44
    45 = 43
    46 = 42
    47 = 41
    48 = 49
    49 = 48
    50 = 38
}

$firstRow = 34
$lastRow = 50

# Only these columns actually vary record-to-record within rows 34-50.
$columns = @("A", "B", "D", "E", "F", "G", "H", "M", "Q", "R", "Z", "AB", "AW", "AX")

# 1) Snapshot the current value of every (column, row) cell before any writes.
$snapshot = @{}
foreach ($col in $columns) {
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $snapshot["$col$r"] = $ws.Range("$col$r").Value2
    }
}

# 2) Write each destination cell from the snapshot of its source row.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    if ($srcRow -ne $destRow) {
        foreach ($col in $columns) {
            $ws.Range("$col$destRow").Value2 = $snapshot["$col$srcRow"]
        }
    }
}
